$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.531682
$ws.Range("H2").Value = 7.595046
$ws.Range("I2").Value = 0.6122191758781785
$ws.Range("J2").Value = 0.6122191758781785
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.970048
$ws.Range("N2").Value = 26.910144
$ws.Range("O2").Value = 0.487108783009476
$ws.Range("P2").Value = 0.4871087830094759
$ws.Range("Q2").Value = 22.709309060736
$ws.Range("R2").Value = 204.383781546624
$ws.Range("S2").Value = 0.2982173376970839
$ws.Range("T2").Value = 0.2982173376970839
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.531682
$ws.Range("H3").Value = 7.595046
$ws.Range("I3").Value = 0.6122191758781785
$ws.Range("J3").Value = 0.6122191758781785
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.012070666666666
$ws.Range("N3").Value = 27.036212
$ws.Range("O3").Value = 0.489390778604016
$ws.Range("P3").Value = 0.489390778604016
$ws.Range("Q3").Value = 22.815697089528
$ws.Range("R3").Value = 205.341273805752
$ws.Range("S3").Value = 0.2996144191593308
$ws.Range("T3").Value = 0.2996144191593308
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.531682
$ws.Range("H4").Value = 7.595046
$ws.Range("I4").Value = 0.6122191758781785
$ws.Range("J4").Value = 0.6122191758781785
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.4327576666666667
$ws.Range("N4").Value = 1.298273
$ws.Range("O4").Value = 0.02350043838650813
$ws.Range("P4").Value = 0.02350043838650813
$ws.Range("Q4").Value = 1.095604795062
$ws.Range("R4").Value = 9.860443155558
$ws.Range("S4").Value = 0.01438741902176392
$ws.Range("T4").Value = 0.01438741902176392
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.603572333333333
$ws.Range("H5").Value = 4.810717
$ws.Range("I5").Value = 0.3877808241218215
$ws.Range("J5").Value = 0.3877808241218215
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.970048
$ws.Range("N5").Value = 26.910144
$ws.Range("O5").Value = 0.487108783009476
$ws.Range("P5").Value = 0.4871087830094759
$ws.Range("Q5").Value = 14.384120801472
$ws.Range("R5").Value = 129.457087213248
$ws.Range("S5").Value = 0.1888914453123921
$ws.Range("T5").Value = 0.1888914453123921
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.603572333333333
$ws.Range("H6").Value = 4.810717
$ws.Range("I6").Value = 0.3877808241218215
$ws.Range("J6").Value = 0.3877808241218215
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 9.012070666666666
$ws.Range("N6").Value = 27.036212
$ws.Range("O6").Value = 0.489390778604016
$ws.Range("P6").Value = 0.489390778604016
$ws.Range("Q6").Value = 14.45150718711156
$ws.Range("R6").Value = 130.063564684004
$ws.Range("S6").Value = 0.1897763594446852
$ws.Range("T6").Value = 0.1897763594446852
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.603572333333333
$ws.Range("H7").Value = 4.810717
$ws.Range("I7").Value = 0.3877808241218215
$ws.Range("J7").Value = 0.3877808241218215
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.4327576666666667
$ws.Range("N7").Value = 1.298273
$ws.Range("O7").Value = 0.02350043838650813
$ws.Range("P7").Value = 0.02350043838650813
$ws.Range("Q7").Value = 0.6939582213045555
$ws.Range("R7").Value = 6.245623991741001
$ws.Range("S7").Value = 0.009113019364744209
$ws.Range("T7").Value = 0.009113019364744209
